$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "56.382.02"
Set-TextValue $ws.Range("E2") "  +3.73%  "
Set-TextValue $ws.Range("D3") "2.501.51"
Set-TextValue $ws.Range("E3") "  +2.69%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "487.59"
Set-TextValue $ws.Range("E5") "  +4.91%  "
Set-TextValue $ws.Range("D6") "145.77"
Set-TextValue $ws.Range("E6") "  +9.97%  "
Set-TextValue $ws.Range("D7") "0.997"
Set-TextValue $ws.Range("E7") "  +0.29%  "
Set-TextValue $ws.Range("D8") "0.510"
Set-TextValue $ws.Range("E8") "  +3.86%  "
Set-TextValue $ws.Range("D9") "2.519.80"
Set-TextValue $ws.Range("E9") "  +2.86%  "
Set-TextValue $ws.Range("D10") "5.68"
Set-TextValue $ws.Range("E10") "  +5.56%  "
Set-TextValue $ws.Range("D11") "0.0974"
Set-TextValue $ws.Range("E11") "  +1.52%  "
Set-TextValue $ws.Range("D12") "0.333"
Set-TextValue $ws.Range("E12") "  +4.10%  "
Set-TextValue $ws.Range("E13") "  +1.18%  "
Set-TextValue $ws.Range("D14") "2.932.89"
Set-TextValue $ws.Range("E14") "  +2.81%  "
Set-TextValue $ws.Range("D15") "56.218.23"
Set-TextValue $ws.Range("E15") "  +3.87%  "
Set-TextValue $ws.Range("D16") "21.11"
Set-TextValue $ws.Range("E16") "  +6.42%  "
Set-TextValue $ws.Range("D17") "0.0000136"
Set-TextValue $ws.Range("E17") "  +2.31%  "
Set-TextValue $ws.Range("D18") "2.513.75"
Set-TextValue $ws.Range("E18") "  +2.75%  "
Set-TextValue $ws.Range("D19") "4.47"
Set-TextValue $ws.Range("E19") "  +6.00%  "
Set-TextValue $ws.Range("D20") "10.28"
Set-TextValue $ws.Range("E20") "  +9.42%  "
Set-TextValue $ws.Range("D21") "320.65"
Set-TextValue $ws.Range("E21") "  +2.10%  "
Set-TextValue $ws.Range("D22") "0.997"
Set-TextValue $ws.Range("E22") "  -0.37%  "
Set-TextValue $ws.Range("D23") "5.81"
Set-TextValue $ws.Range("E23") "  +7.71%  "
Set-TextValue $ws.Range("D24") "58.63"
Set-TextValue $ws.Range("E24") "  +3.10%  "
Set-TextValue $ws.Range("D25") "0.411"
Set-TextValue $ws.Range("E25") "  +6.77%  "
Set-TextValue $ws.Range("E26") "  +7.83%  "
Set-TextValue $ws.Range("D27") "0.997"
Set-TextValue $ws.Range("E27") "  -1.04%  "
Set-TextValue $ws.Range("D28") "2.607.71"
Set-TextValue $ws.Range("E28") "  +3.65%  "
Set-TextValue $ws.Range("D29") "7.53"
Set-TextValue $ws.Range("E29") "  +4.74%  "
Set-TextValue $ws.Range("D30") "0.0₃0787"
Set-TextValue $ws.Range("E30") "  +7.92%  "
Set-TextValue $ws.Range("E31") "  +0.34%  "
Set-TextValue $ws.Range("D32") "148.25"
Set-TextValue $ws.Range("E32") "  -1.70%  "
Set-TextValue $ws.Range("D33") "18.34"
Set-TextValue $ws.Range("E33") "  +3.63%  "
Set-TextValue $ws.Range("E34") "  +7.18%  "
Set-TextValue $ws.Range("D35") "5.21"
Set-TextValue $ws.Range("E35") "  +2.96%  "
Set-TextValue $ws.Range("E36") "  +8.46%  "
Set-TextValue $ws.Range("E37") "  +4.68%  "
Set-TextValue $ws.Range("D38") "0.868"
Set-TextValue $ws.Range("E38") "  +8.28%  "
Set-TextValue $ws.Range("D39") "34.15"
Set-TextValue $ws.Range("E39") "  +1.66%  "
Set-TextValue $ws.Range("E40") "  +7.35%  "
Set-TextValue $ws.Range("D41") "0.618"
Set-TextValue $ws.Range("E41") "  +1.58%  "
Set-TextValue $ws.Range("D42") "0.996"
Set-TextValue $ws.Range("E42") "  +0.46%  "
Set-TextValue $ws.Range("E43") "  +4.73%  "
Set-TextValue $ws.Range("E44") "  +6.35%  "
Set-TextValue $ws.Range("E45") "  +10.12%  "
Set-TextValue $ws.Range("D46") "263.54"
Set-TextValue $ws.Range("E46") "  +20.77%  "
Set-TextValue $ws.Range("D49") "0.0908"
Set-TextValue $ws.Range("E49") "  +4.31%  "
Set-TextValue $ws.Range("D50") "1.922.65"
Set-TextValue $ws.Range("E50") "  -2.73%  "
Set-TextValue $ws.Range("D51") "17.65"
Set-TextValue $ws.Range("E51") "  +6.35%  "

# Rows 47/48: WhiteBITCoin and VeChain swapped positions with updated data
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0228"
Set-TextValue $ws.Range("E47") "  +3.19%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D48") "10.16"
Set-TextValue $ws.Range("E48") "  -0.42%  "
